$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Fix a typo in the generated input-file names referenced in columns D (dbExcel)
# and E (WebExcel): "..._Genomi_..." -> "..._Genomic_..." for rows 2-6.
for ($row = 2; $row -le 6; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)

    $dCell.Value = $dCell.Value2.Replace("Genomi_", "Genomic_")
    $eCell.Value = $eCell.Value2.Replace("Genomi_", "Genomic_")
}

# Update the saved view state of the sheet (scrolled/selected cell).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D16").Select()
